$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z")

# ---- Row 56: 2025-09-28, 四方坪站充电量(kw) ----
$ws.Range("A56").Value = 45928
$ws.Range("A56").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B56").Value = "四方坪站充电量(kw)"

$row56 = @(833.23699999999997,911.19700000000023,534.726,426.75,227.999,639.88800000000015,437.94,187.88300000000004,132.74799999999999,254.53000000000003,123.41999999999999,179.68700000000001,575.58800000000019,1449.7339999999992,396.62599999999998,357.48900000000003,305.55799999999999,290.072,171.59099999999998,141.47599999999997,60.64,191.029,170.19,85.036000000000001)

$ws.Range("C56:Z56").NumberFormat = "0.00"
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])56").Value = $row56[$i]
}

# ---- Row 57: 2025-09-28, 高岭站充电量(kw) ----
$ws.Range("A57").Value = 45928
$ws.Range("A57").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B57").Value = "高岭站充电量(kw)"

$row57 = @(653.05199999999991,590.49400000000003,310.74300000000005,25.172000000000001,87.162999999999997,209.50799999999998,104.65199999999999,298.06700000000001,342.60500000000002,199.321,102.791,163.52699999999999,315.39200000000005,667.99399999999991,171.29500000000002,451.23099999999994,478.75299999999999,157.30799999999999,216.98599999999999,45.960999999999999,100.70899999999999,36.061999999999998,142.423,0)

$ws.Range("C57:Z57").NumberFormat = "0.00"
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])57").Value = $row57[$i]
}

# Matches the author's final selection after entering the new rows.
$ws.Range("B61").Select()
